$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# Sheet 1 updates
$ws1.Range("F2").Value = 93
$ws1.Range("F3").Value = 72
$ws1.Range("F4").Value = 969
$ws1.Range("F5").Value = 1272
$ws1.Range("F6").Value = 1757
$ws1.Range("F7").Value = 931
$ws1.Range("F8").Value = 578
$ws1.Range("F9").Value = 2649
$ws1.Range("F10").Value = 746
$ws1.Range("F12").Value = 574
$ws1.Range("F13").Value = 39
$ws1.Range("F15").Value = 345
$ws1.Range("F16").Value = 340
$ws1.Range("F18").Value = 2133
$ws1.Range("F21").Value = 8
$ws1.Range("F22").Value = 2629
$ws1.Range("F25").Value = 27
$ws1.Range("F26").Value = 522
$ws1.Range("F28").Value = 469
$ws1.Range("F31").Value = 530
$ws1.Range("F32").Value = 552
$ws1.Range("F33").Value = 215
$ws1.Range("F35").Value = 347
$ws1.Range("F36").Value = 4601
$ws1.Range("F37").Value = 172

# Sheet 2 updates
$ws2.Range("F2").Value = 393
$ws2.Range("F4").Value = 4211
$ws2.Range("F12").Value = 16
$ws2.Range("F13").Value = 320
$ws2.Range("F14").Value = 330
$ws2.Range("F17").Value = 5
$ws2.Range("F21").Value = 5
$ws2.Range("F26").Value = 1
$ws2.Range("F28").Value = 16
$ws2.Range("F29").Value = 278

# Sheet 3 updates
$ws3.Range("F4").Value = 1433
$ws3.Range("F6").Value = 539
$ws3.Range("F7").Value = 149

# Sheet 4 updates
$ws4.Range("F2").Value = 1433
$ws4.Range("F3").Value = 539
$ws4.Range("F4").Value = 393
$ws4.Range("F6").Value = 93
$ws4.Range("F8").Value = 72
$ws4.Range("F9").Value = 969
$ws4.Range("F10").Value = 1272
$ws4.Range("F11").Value = 1757
$ws4.Range("F14").Value = 149
$ws4.Range("F15").Value = 931
$ws4.Range("F16").Value = 578
$ws4.Range("F17").Value = 2649
$ws4.Range("F18").Value = 746
$ws4.Range("F19").Value = 574
$ws4.Range("F20").Value = 345
$ws4.Range("F21").Value = 340
$ws4.Range("F22").Value = 320
$ws4.Range("F23").Value = 330
$ws4.Range("F25").Value = 2133
$ws4.Range("F27").Value = 5
$ws4.Range("F28").Value = 8
$ws4.Range("F30").Value = 2629
$ws4.Range("F38").Value = 469
$ws4.Range("F40").Value = 16
$ws4.Range("F41").Value = 278
$ws4.Range("F43").Value = 215
$ws4.Range("F44").Value = 347
$ws4.Range("F45").Value = 4601
$ws4.Range("F46").Value = 172
